# other_species.xlsx update:
#  - rename the "form_id" setting to "table_id" on the settings sheet
#  - add a new "properties" sheet (after "settings") that will drive
#    properties.csv generation, with a minimal partition/aspect/key/type/value
#    table describing the Table's default column order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. settings sheet: rename the form_id setting to table_id
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "table_id"
[void]$settings.Range("A3").Select()

# ---------------------------------------------------------------------------
# 2. add the new "properties" sheet, placed after "settings" (last tab)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$properties = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$properties.Name = "properties"

# Header row
$properties.Range("A1").Value = "partition"
$properties.Range("B1").Value = "aspect"
$properties.Range("C1").Value = "key"
$properties.Range("D1").Value = "type"
$properties.Range("E1").Value = "value"

# Data row: the default colOrder property for the Table partition
$properties.Range("A2").Value = "Table"
$properties.Range("B2").Value = "default"
$properties.Range("C2").Value = "colOrder"
$properties.Range("D2").Value = "array"
$properties.Range("E2").Value = '["OS_FOL_date","OS_FOL_B_focal_AnimID","OS_time_begin","OS_time_end","OS_OSL_local_species_name","OS_local_species_name_written","OS_duration","OS_comments"]'

[void]$properties.Range("E3").Select()
[void]$properties.Activate()
